# Update the date line (unique text in the document body, outside the table).
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-11-11 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-11-12 Sunday", 2)

# Update the division problems inside the table by directly addressing each
# cell (row, column). This avoids ambiguity from duplicate cell text values
# that appear as both an old value in one cell and a new value in another.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "19÷6=3, 1" },
    @{ Row = 1;  Col = 2; Text = "88÷2=44, 0" },
    @{ Row = 1;  Col = 3; Text = "39÷6=6, 3" },
    @{ Row = 1;  Col = 4; Text = "61÷9=6, 7" },
    @{ Row = 1;  Col = 5; Text = "68÷5=13, 3" },

    @{ Row = 5;  Col = 1; Text = "17÷9=1, 8" },
    @{ Row = 5;  Col = 2; Text = "33÷2=16, 1" },
    @{ Row = 5;  Col = 3; Text = "18÷6=3, 0" },
    @{ Row = 5;  Col = 4; Text = "31÷2=15, 1" },
    @{ Row = 5;  Col = 5; Text = "74÷9=8, 2" },

    @{ Row = 9;  Col = 1; Text = "68÷8=8, 4" },
    @{ Row = 9;  Col = 2; Text = "33÷3=11, 0" },
    @{ Row = 9;  Col = 3; Text = "51÷4=12, 3" },
    @{ Row = 9;  Col = 4; Text = "58÷5=11, 3" },
    @{ Row = 9;  Col = 5; Text = "13÷2=6, 1" },

    @{ Row = 13; Col = 1; Text = "87÷4=21, 3" },
    @{ Row = 13; Col = 2; Text = "53÷9=5, 8" },
    @{ Row = 13; Col = 3; Text = "82÷2=41, 0" },
    @{ Row = 13; Col = 4; Text = "18÷9=2, 0" },
    @{ Row = 13; Col = 5; Text = "45÷7=6, 3" },

    @{ Row = 17; Col = 1; Text = "96÷8=12, 0" },
    @{ Row = 17; Col = 2; Text = "19÷7=2, 5" },
    @{ Row = 17; Col = 3; Text = "71÷3=23, 2" },
    @{ Row = 17; Col = 4; Text = "88÷5=17, 3" },
    @{ Row = 17; Col = 5; Text = "41÷8=5, 1" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
